$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.Style = "Normal"
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "43.681.86"
Set-TextValue "E2" "  +3.97%  "
Set-TextValue "D3" "2.229.53"
Set-TextValue "E3" "  +3.12%  "
Set-TextValue "E4" "  -0.26%  "
Set-TextValue "D5" "259.40"
Set-TextValue "E5" "  +2.72%  "
Set-TextValue "D6" "81.24"
Set-TextValue "E6" "  +11.20%  "
Set-TextValue "E7" "  +2.94%  "
Set-TextValue "E9" "  +2.95%  "
Set-TextValue "D10" "43.40"
Set-TextValue "E10" "  +9.17%  "
Set-TextValue "D11" "0.0926"
Set-TextValue "E11" "  +1.97%  "
Set-TextValue "E12" "  +4.30%  "
Set-TextValue "E13" "  +2.66%  "
Set-TextValue "D14" "2.561.57"
Set-TextValue "E14" "  +2.75%  "
Set-TextValue "D15" "14.64"
Set-TextValue "E15" "  +3.09%  "
Set-TextValue "D16" "2.219.34"
Set-TextValue "E16" "  +2.03%  "
Set-TextValue "D17" "0.786"
Set-TextValue "E17" "  +2.35%  "
Set-TextValue "D18" "43.609.10"
Set-TextValue "E18" "  +4.00%  "
Set-TextValue "E19" "  +2.45%  "
Set-TextValue "D20" "71.12"
Set-TextValue "E20" "  +0.67%  "
Set-TextValue "E21" "  +3.45%  "
Set-TextValue "E22" "  +9.93%  "
Set-TextValue "D23" "232.56"
Set-TextValue "E23" "  +2.82%  "
Set-TextValue "D24" "9.26"
Set-TextValue "E24" "  -3.04%  "
Set-TextValue "E25" "  +0.14%  "
Set-TextValue "D26" "10.79"
Set-TextValue "E26" "  +2.71%  "
Set-TextValue "D27" "41.35"
Set-TextValue "E27" "  +12.34%  "
Set-TextValue "E28" "  +1.36%  "
Set-TextValue "D29" "2.24"
Set-TextValue "E29" "  +2.48%  "
Set-TextValue "E30" "  -0.25%  "
Set-TextValue "D31" "172.60"
Set-TextValue "E31" "  +2.33%  "
Set-TextValue "D32" "0.0897"
Set-TextValue "E32" "  +12.46%  "
Set-TextValue "D33" "20.61"
Set-TextValue "E33" "  +3.53%  "
Set-TextValue "D34" "5.31"
Set-TextValue "E34" "  +4.24%  "
Set-TextValue "D35" "0.115"
Set-TextValue "E35" "  +8.12%  "
Set-TextValue "B36" "VeChain"
Set-TextValue "C36" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D36" "0.0371"
Set-TextValue "E36" "  +13.71%  "
Set-TextValue "B37" "Stellar"
Set-TextValue "C37" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D37" "0.123"
Set-TextValue "E37" "  +2.47%  "
Set-TextValue "D38" "4.59"
Set-TextValue "E38" "  +7.71%  "
Set-TextValue "B39" "Celestia"
Set-TextValue "C39" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D39" "12.94"
Set-TextValue "E39" "  +7.77%  "
Set-TextValue "B40" "NEARProtocol"
Set-TextValue "C40" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D40" "2.98"
Set-TextValue "E40" "  +25.29%  "
Set-TextValue "D41" "2.13"
Set-TextValue "E41" "  +3.66%  "
Set-TextValue "D42" "63.44"
Set-TextValue "E42" "  +8.08%  "
Set-TextValue "D43" "5.51"
Set-TextValue "E43" "  +7.40%  "
Set-TextValue "E44" "  +3.26%  "
Set-TextValue "D45" "103.98"
Set-TextValue "E45" "  +3.21%  "
Set-TextValue "E46" "  +2.30%  "
Set-TextValue "E47" "  +2.23%  "
Set-TextValue "D48" "1.13"
Set-TextValue "E48" "  +4.00%  "
Set-TextValue "B49" "Stacks"
Set-TextValue "C49" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D49" "1.56"
Set-TextValue "E49" "  +27.92%  "
Set-TextValue "B50" "WOONetwork"
Set-TextValue "C50" "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextValue "D50" "0.443"
Set-TextValue "E50" "  -5.73%  "
Set-TextValue "E51" "  +3.49%  "
